{"js": "context.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\nWrite-Output $d.Paragraphs.Count\n"}
